# Update column G ("K" = strikeouts) values for rows 2-71 on Sheet1.
# The workbook's save_data generation was regenerated to compute K from
# the pitch-by-pitch data (s_vals) instead of a raw "Strike#" field, which
# changes most of the previously stored values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 1
    3  = 0
    4  = 0
    5  = 1
    6  = 3
    7  = 2
    8  = 2
    9  = 1
    10 = 2
    11 = 0
    13 = 1
    14 = 0
    15 = 0
    16 = 1
    17 = 0
    18 = 1
    19 = 1
    20 = 1
    21 = 0
    22 = 2
    23 = 1
    24 = 1
    25 = 0
    26 = 1
    27 = 0
    28 = 1
    29 = 1
    30 = 1
    31 = 0
    32 = 2
    33 = 1
    34 = 1
    35 = 1
    37 = 1
    38 = 1
    39 = 1
    40 = 2
    42 = 1
    43 = 0
    44 = 0
    45 = 3
    46 = 0
    47 = 2
    48 = 0
    49 = 0
    50 = 0
    51 = 3
    52 = 1
    53 = 2
    54 = 0
    55 = 1
    56 = 1
    57 = 1
    58 = 2
    59 = 0
    60 = 2
    61 = 2
    62 = 1
    63 = 1
    64 = 1
    65 = 0
    66 = 1
    67 = 0
    68 = 2
    69 = 1
    70 = 1
    71 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
